# "Debut PCB. Groupement et placement des composantes."
# Two resistors in the BOM are re-specified from 1206(3216) package parts
# to 0402(1005) package parts (PCB layout/placement pass):
#   - "res 2k"  (row 27): yageo RC1206FR-072KL  -> panasonic ERJ-2RKF2001X
#   - "res 1k5" (row 29): panasonic ERA-8AEB152V -> panasonic ERA-2AEB152X
#     (also a unit price drop from 0.39 to 0.18 for the new part)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 first, so its new shared string ("ERA-2AEB152X...") is interned
# ahead of row 27's ("ERJ-2RKF2001X...") - matches the order newly added
# strings appear in the saved workbook.
$ws.Range("D29").Value = 0.18
$ws.Range("F29").Value = "https://www.digikey.ca/en/products/detail/panasonic-electronic-components/ERA-2AEB152X/1706009"
$ws.Range("H29").Value = "0402(1005)"

$ws.Range("F27").Value = "https://www.digikey.ca/en/products/detail/panasonic-electronic-components/ERJ-2RKF2001X/192194"
$ws.Range("H27").Value = "0402(1005)"

# Leave the selection where the author's work ended up.
$ws.Range("I29").Select()
